$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1185646666666667
$ws.Range("H2").Value = 0.355694
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.164505666666667
$ws.Range("N2").Value = 6.493517000000001
$ws.Range("O2").Value = 0.5225358117353504
$ws.Range("P2").Value = 0.5225358117353504
$ws.Range("Q2").Value = 0.2566338928664444
$ws.Range("R2").Value = 2.309705035798
$ws.Range("S2").Value = 0.5225358117353504
$ws.Range("T2").Value = 0.5225358117353504

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1185646666666667
$ws.Range("H3").Value = 0.355694
$ws.Range("O3").Value = 0.2046934834760502
$ws.Range("P3").Value = 0.2046934834760502
$ws.Range("Q3").Value = 0.1005314551253333
$ws.Range("R3").Value = 0.9047830961279998
$ws.Range("S3").Value = 0.2046934834760502
$ws.Range("T3").Value = 0.2046934834760502

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1185646666666667
$ws.Range("H4").Value = 0.355694
$ws.Range("M4").Value = 1.129901
$ws.Range("N4").Value = 3.389703
$ws.Range("O4").Value = 0.2727707047885994
$ws.Range("P4").Value = 0.2727707047885994
$ws.Range("Q4").Value = 0.1339663354313333
$ws.Range("R4").Value = 1.205697018882
$ws.Range("S4").Value = 0.2727707047885994
$ws.Range("T4").Value = 0.2727707047885994
